# Change default "IAEA-C1" sample ID entries to "IAEA-C2" across the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$targetCells = @("A4", "F10", "A13", "F21", "A27", "F32", "A40", "A50")
foreach ($cellRef in $targetCells) {
    $cell = $ws.Range($cellRef)
    if ($cell.Value2 -eq "IAEA-C1") {
        $cell.Value = "IAEA-C2"
    }
}

# Update the active selection to reflect the new cursor position (F10).
$ws.Range("F10").Select()
